# Update the summary table on the active sheet with the latest
# nest-measure statistics (new season-level latency/relative-timing
# aggregates). The table layout (Description + Winter/Spring/Summer/
# Autumn/All columns, Latency-to-AB / AB-Rel-to-Sunrise / Latency-to-QB /
# QB-Rel-to-Sunset rows) is unchanged -- only the numeric cell contents
# are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Description", "Winter", "Spring", "Summer", "Autumn", "All"),
    @("Latency to AB", "NaN ± NaN", "47.21 ± 30.09", "33.79 ± 34.87", "30.55 ± 12.44", "37.25 ± 31.27"),
    @("AB Rel. to Sunrise", "NaN ± NaN", "58.72 ± 68.03", "60.91 ± 68.89", "-21.12 ± 18.70", "45.69 ± 69.94"),
    @("Latency to QB", "NaN ± NaN", "90.17 ± 113.51", "80.50 ± 92.15", "41.82 ± 33.58", "76.54 ± 93.67"),
    @("QB Rel. to Sunset", "NaN ± NaN", "-137.56 ± 96.61", "-81.92 ± 76.38", "-10.62 ± 52.95", "-85.98 ± 90.32")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws.Columns.Item(3).ColumnWidth = 13.140625
